$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row with ticker value "GRT-USD" at A69
$ws.Range("A69").Value = "GRT-USD"
